# Fix: DCASP fechada para a PM.
# Appends 16 new expense records (rows 383-398) to the "ConsorcioDespesas"
# table on the "Despesas" sheet, extending it from A1:K382 to A1:K398.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Despesas")
$tbl = $ws.ListObjects.Item("ConsorcioDespesas")

# consorcio, data_base, funcao, subfuncao, ndo, empenhado, liquidado, pago
$data = @(
    ,@("COFRON", 45657, 4, 122, 319011010100, 901.39, 901.39, 901.39)
    ,@("COFRON", 45657, 4, 122, 319013010100, 108.4, 108.4, 108.4)
    ,@("COFRON", 45657, 4, 122, 319013020100, 378.58, 378.58, 378.58)
    ,@("COFRON", 45657, 4, 122, 339014140000, 0, 0, 0)
    ,@("COFRON", 45657, 4, 122, 339030000000, 0, 0, 0)
    ,@("COFRON", 45657, 4, 122, 339033010000, 0, 0, 0)
    ,@("COFRON", 45657, 4, 122, 339039000000, 3.61, 88.69, 130.4)
    ,@("COFRON", 45657, 4, 122, 339039990100, 59.11, 59.11, 59.11)
    ,@("COFRON", 45657, 4, 122, 339039400000, 0, 286.32, 166.83)
    ,@("COFRON", 45657, 4, 122, 339046010100, 53.48, 53.48, 53.48)
    ,@("COFRON", 45657, 4, 122, 339047000000, 0, 0, 0)
    ,@("COFRON", 45657, 4, 122, 339049010000, 1.96, 1.96, 1.96)
    ,@("COFRON", 45657, 4, 122, 449052000000, 190.74, 190.74, 190.74)
    ,@("COFRON", 45657, 10, 302, 334041390500, 9273.68, 9273.68, 9273.68)
    ,@("COFRON", 45657, 10, 302, 334041391100, 865.54, 865.54, 865.54)
    ,@("COFRON", 45657, 10, 302, 334041391000, 655.91, 655.91, 655.91)
)

$lastDataRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1

foreach ($rec in $data) {

    $listRow = $tbl.ListRows.Add()
    $rowRange = $listRow.Range
    $rowNum = $rowRange.Row

    # Copy number formats (date, ndo, currency) from the previous data row
    # so the new row keeps the same visual formatting as the rest of the table.
    $ws.Range("B$lastDataRow").Copy() | Out-Null
    $ws.Range("B$rowNum").PasteSpecial(-4122) | Out-Null

    $ws.Range("E$lastDataRow").Copy() | Out-Null
    $ws.Range("E$rowNum").PasteSpecial(-4122) | Out-Null

    $ws.Range("F$($lastDataRow):H$lastDataRow").Copy() | Out-Null
    $ws.Range("F$($rowNum):H$rowNum").PasteSpecial(-4122) | Out-Null

    $excel.CutCopyMode = 0

    # Values
    $rowRange.Cells.Item(1, 1).Value = $rec[0]
    $rowRange.Cells.Item(1, 2).Value = $rec[1]
    $rowRange.Cells.Item(1, 3).Value = $rec[2]
    $rowRange.Cells.Item(1, 4).Value = $rec[3]
    $rowRange.Cells.Item(1, 5).Value = $rec[4]
    $rowRange.Cells.Item(1, 6).Value = $rec[5]
    $rowRange.Cells.Item(1, 7).Value = $rec[6]
    $rowRange.Cells.Item(1, 8).Value = $rec[7]

    # Calculated columns (ano / bimestre / mes)
    $iCell = $rowRange.Cells.Item(1, 9)
    $jCell = $rowRange.Cells.Item(1, 10)
    $kCell = $rowRange.Cells.Item(1, 11)

    $iCell.Formula = "=YEAR(ConsorcioDespesas[[#This Row],[data_base]])"
    $jCell.Formula = "=_xlfn.SWITCH(MONTH(ConsorcioDespesas[[#This Row],[data_base]]),1,1,2,1,3,2,4,2,5,3,6,3,7,4,8,4,9,5,10,5,11,6,12,6)"
    $kCell.Formula = "=MONTH(ConsorcioDespesas[[#This Row],[data_base]])"

    # Re-create the "Virgula" calculated-column style used by the rest of the
    # table (General number format, default font, linked to the Virgula
    # cell style) for the ano/bimestre/mes columns.
    foreach ($c in @($iCell, $jCell, $kCell)) {
        $c.NumberFormat = "General"
        $c.Style = "Vírgula"
    }

    $lastDataRow = $rowNum
}

$excel.CutCopyMode = 0

# Match the sheet view state recorded in the edited workbook.
$ws.Application.ActiveWindow.ScrollRow = 363
$ws.Range("F396").Select()
